$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07844425225260536
$ws.Range("C2").Value = 0.0782501778507621
$ws.Range("B3").Value = 23.13902426738357
$ws.Range("C3").Value = 23.13854921214313
$ws.Range("B4").Value = 151.7000087709266
$ws.Range("C4").Value = 151.7000486805258
$ws.Range("B5").Value = 0.1859816921622611
$ws.Range("C5").Value = 0.1857631005249513
$ws.Range("B6").Value = 2.549980640079351
$ws.Range("C6").Value = 2.550019590994558
$ws.Range("B7").Value = 0.7692259974068437
$ws.Range("C7").Value = 0.7631674057556723
$ws.Range("B9").Value = 1.825554439925985
$ws.Range("C9").Value = 1.824893393937746
$ws.Range("B11").Value = 907.4136508152212
$ws.Range("C11").Value = 910.2021873357731
$ws.Range("B12").Value = 0.9749011072409366
$ws.Range("C12").Value = 1.075446053817627
$ws.Range("B13").Value = 0.9402340008544148
$ws.Range("C13").Value = 0.9951654437442093
$ws.Range("B14").Value = 2.438056001776486
$ws.Range("C14").Value = 2.438471660727712
$ws.Range("B15").Value = 0.9142119329940921
$ws.Range("C15").Value = 0.914256720186948
$ws.Range("B16").Value = 0.07815659036814897
$ws.Range("C16").Value = 0.07882760274675904
$ws.Range("B19").Value = 0.4549029917249756
$ws.Range("C19").Value = 0.4536924731850582
$ws.Range("B20").Value = 0.1259456542906707
$ws.Range("C20").Value = 0.1269705784710714
$ws.Range("B21").Value = 0.1288642172820511
$ws.Range("C21").Value = 0.1252681519406583
$ws.Range("B22").Value = 5.109687419198908
$ws.Range("C22").Value = 5.10790859792463
$ws.Range("B23").Value = -0.004097353059522491
$ws.Range("C23").Value = 0.00804272403612475
$ws.Range("B24").Value = 0.4033872098114001
$ws.Range("C24").Value = 0.4056507790300997
$ws.Range("B25").Value = 28.30654607309969
$ws.Range("C25").Value = 28.30365333274037
$ws.Range("B26").Value = 29.29958102108359
$ws.Range("C26").Value = 29.29978708603336
$ws.Range("B27").Value = 0.08189844844316098
$ws.Range("C27").Value = 0.08213406191251332
$ws.Range("B28").Value = 0.4546600195179318
$ws.Range("C28").Value = 0.4544837018245149
$ws.Range("B29").Value = 1.70654392870338
$ws.Range("C29").Value = 1.706871611059007
$ws.Range("B30").Value = 3.598086841575389
$ws.Range("C30").Value = 3.596349917560963
$ws.Range("B31").Value = 12.06547291340627
$ws.Range("C31").Value = 12.0654267495494
$ws.Range("B32").Value = 33.39448477844254
$ws.Range("C32").Value = 33.37300287709702
$ws.Range("B33").Value = 73843.24066627219
$ws.Range("C33").Value = 73797.24431837903
$ws.Range("B34").Value = 9.004229623816084
$ws.Range("C34").Value = 8.999396750164943
$ws.Range("B35").Value = 87.55305656026985
$ws.Range("C35").Value = 87.95448089355044
$ws.Range("B36").Value = 124.2173221769794
$ws.Range("C36").Value = 122.8319242020687
